$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their original text formatting
# (values like "69.637.05" or "3.30" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '69.637.05'
$ws.Range('E2').Value = '  +2.08%  '
$ws.Range('D3').Value = '3.377.52'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '580.67'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '179.19'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.594'
$ws.Range('E8').Value = '  +0.74%  '
$ws.Range('D9').Value = '0.199'
$ws.Range('E9').Value = '  +8.48%  '
$ws.Range('D10').Value = '0.589'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').Value = '48.32'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '0.0000284'
$ws.Range('E12').Value = '  +4.01%  '
$ws.Range('D13').Value = '685.48'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').Value = '8.62'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').Value = '3.921.39'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').Value = '69.584.38'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').Value = '0.120'
$ws.Range('D18').Value = '3.387.46'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('D19').Value = '17.69'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = '11.25'
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('D21').Value = '0.910'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '17.25'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').Value = '101.67'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').Value = '3.89'
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('D26').Value = '2.70'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '9.70'
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('D28').Value = '33.57'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').Value = '8.74'
$ws.Range('E29').Value = '  +2.67%  '
$ws.Range('D30').Value = '6.90'
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('D31').Value = '3.84'
$ws.Range('E31').Value = '  +16.81%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '11.05'
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '554.07'
$ws.Range('E33').Value = '  -1.90%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = '57.82'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '3.601.28'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('D39').Value = '35.32'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('D40').Value = '0.0₃0728'
$ws.Range('E40').Value = '  +8.34%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '3.30'
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.74'
$ws.Range('E42').Value = '  +4.63%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = '3.37'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0425'
$ws.Range('E44').Value = '  +2.83%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.336'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '2.66'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.129'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '1.38'
$ws.Range('E48').Value = '  +3.77%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '130.05'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = '2.57'
$ws.Range('E51').Value = '  -0.23%  '
